$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) values of the columns that change,
# for rows 2, 3 and 4, before overwriting anything.
$cols = @("D", "J", "K", "L", "M", "O", "P")

$row2 = @{}
$row3 = @{}
$row4 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range("${col}2").Value2
    $row3[$col] = $ws.Range("${col}3").Value2
    $row4[$col] = $ws.Range("${col}4").Value2
}

# Cyclic rotation: new row2 = old row3, new row3 = old row4, new row4 = old row2
foreach ($col in $cols) {
    $ws.Range("${col}2").Value2 = $row3[$col]
    $ws.Range("${col}3").Value2 = $row4[$col]
    $ws.Range("${col}4").Value2 = $row2[$col]
}
